$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 4 new rows of accelerometer samples above the existing data
# (before the old row 2), pushing the old rows 2-21 down to rows 6-25.
$ws.Rows.Item(2).Resize(4).Insert()
$ws.Rows.Item(2).Resize(4).ClearFormats()

$newTop = @(
    @(0.08896994590759319, -0.11671480536461, 0.04848458990454618),
    @(-0.06001234054565427, 0.2669965513050556, 0.1866837395355105),
    @(-0.111260414123535, 0.1850093007087709, -0.02702043950557642),
    @(0.04686117172241161, -0.03032520040869657, 0.02842492796480637)
)

for ($i = 0; $i -lt $newTop.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $newTop[$i][0]
    $ws.Cells.Item($r, 2).Value = $newTop[$i][1]
    $ws.Cells.Item($r, 3).Value = $newTop[$i][2]
}

# --- Append 6 more new rows of samples after the (now shifted) last row,
# i.e. rows 26-31.
$newBottom = @(
    @(-0.5248832702637101, -2.431172959506501, -2.120043188333502),
    @(-0.4282075166702048, -0.7951091900468019, -1.767483308911331),
    @(-2.543609619140639, 0.4462372660637008, -1.760738492012012),
    @(0.2467263936996389, -0.1943315342068692, -2.037849001586441),
    @(0.9171624183654843, 0.3103487230837345, 0.1225722581148094),
    @(-0.926007807254792, 0.2953229788690807, 2.136403992772098)
)

$startRow = 26
for ($i = 0; $i -lt $newBottom.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newBottom[$i][0]
    $ws.Cells.Item($r, 2).Value = $newBottom[$i][1]
    $ws.Cells.Item($r, 3).Value = $newBottom[$i][2]
}
